# CAI_Workout_Regimen.xlsx - reorder exercises within each day
# (Day 3: move "Dips" from position 2 to position 5)
#
# The author's edit also padded the table with a leading blank row/column,
# which we reproduce by inserting a new row above row 1 and a new column
# to the left of column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workout Plan")

# --- Add a spacer row above and a spacer column to the left of the table ---
$ws.Rows("1:1").Insert()
$ws.Columns("A:A").Insert()
$ws.Columns("A:A").ColumnWidth = 3.6666666666666665

# --- Reorder Day 3 exercises: move "Dips" (row 20) to sit right after
#     "Barbell Back Extension" (row 23), i.e. become the 5th exercise of the day ---
$ws.Range("B24:H24").Insert()
$ws.Range("B20:H20").Cut($ws.Range("B24:H24"))
$ws.Range("B20:H20").Delete()

# Renumber the "Order" column for the rows that shifted within Day 3
$ws.Range("C20").Value = 2
$ws.Range("C21").Value = 3
$ws.Range("C22").Value = 4
$ws.Range("C23").Value = 5

# --- Restore cursor/selection position like the source workbook ---
$null = $ws.Range("I36").Select()
